$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.608.59"
$ws.Range("E2").Value = "  -0.25%  "

# Row 3
$ws.Range("D3").Value = "3.144.49"
$ws.Range("E3").Value = "  -0.34%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").Value = "3.142.85"
$ws.Range("E8").Value = "  -0.38%  "

# Row 9
$ws.Range("E9").Value = "  -0.45%  "

# Row 10
$ws.Range("E10").Value = "  -2.73%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.38%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.499"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.97%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.46%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.15%  "

# Row 15
$ws.Range("D15").Value = "3.661.67"
$ws.Range("E15").Value = "  -0.16%  "

# Row 16
$ws.Range("D16").Value = "64.696.17"

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.144.13"
$ws.Range("E17").Value = "  -0.52%  "

# Row 18
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.75%  "

# Row 19
$ws.Range("E19").Value = "  +0.31%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "502.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.86%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.711"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.51%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.40%  "

# Row 26
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.08%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.35%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.16%  "

# Row 30
$ws.Range("E30").Value = "  +4.96%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.09%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.04%  "

# Row 33
$ws.Range("E33").Value = "  +0.76%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.70%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.15%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0888"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.66%  "

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "477.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.43%  "

# Row 39
$ws.Range("E39").Value = "  -2.19%  "

# Row 40
$ws.Range("E40").Value = "  -3.47%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "

# Row 42
$ws.Range("D42").Value = "2.997.92"
$ws.Range("E42").Value = "  -3.74%  "

# Row 43
$ws.Range("E43").Value = "  -4.52%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.280"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.75%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "27.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.35%  "

# Row 47
$ws.Range("D47").Value = "0.0₃0578"
$ws.Range("E47").Value = "  -0.39%  "

# Row 49
$ws.Range("E49").Value = "  -1.81%  "

# Row 50
$ws.Range("E50").Value = "  -3.10%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.62%  "
